# fix(gui) step 1 and 2
# - bump the date in A1 by one day
# - update the price column (D) for the "ARANDELAS PLANAS" and
#   "ARANDELA CHAPISTA" price lists (rows 33-47 and 51-54)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: date bump
$ws.Range("A1").Value = 45309

# Step 2: price updates
$ws.Cells.Item(33, 4).Value = 5030.479
$ws.Cells.Item(34, 4).Value = 3779.603
$ws.Cells.Item(35, 4).Value = 3261.594
$ws.Cells.Item(36, 4).Value = 2929.043
$ws.Cells.Item(37, 4).Value = 2929.043
$ws.Cells.Item(38, 4).Value = 2470.491
$ws.Cells.Item(39, 4).Value = 2470.491
$ws.Cells.Item(40, 4).Value = 2470.491
$ws.Cells.Item(41, 4).Value = 2470.491
$ws.Cells.Item(42, 4).Value = 2470.491
$ws.Cells.Item(43, 4).Value = 2470.491
$ws.Cells.Item(44, 4).Value = 2470.491
$ws.Cells.Item(45, 4).Value = 2948.224
$ws.Cells.Item(46, 4).Value = 2948.224
$ws.Cells.Item(47, 4).Value = 2948.224

$ws.Cells.Item(51, 4).Value = 4297.629
$ws.Cells.Item(52, 4).Value = 3907.523
$ws.Cells.Item(53, 4).Value = 3261.594
$ws.Cells.Item(54, 4).Value = 3261.594
